$d = $word.ActiveDocument
$nbsp = [char]0x00A0

# ----------------------------------------------------------------------
# Locate the paragraph that reads:
#   "Product Owner<nbsp>: Prof. Dr. Christian Kruse"
# (two runs: "Product Owner<nbsp>: " and "Prof. Dr. Christian Kruse")
# and rewrite it to:
#   "Produkt Owner<nbsp>: Herr Prof. Dr. Christian Kruse"
# split across five runs:
#   "Produkt Owner<nbsp>: Herr" | " Prof. Dr." | " " | "Christian " | "Kruse"
# ----------------------------------------------------------------------

$oldRun1 = "Product Owner" + $nbsp + ": "
$oldRun2 = "Prof. Dr. Christian Kruse"
$oldFull = $oldRun1 + $oldRun2

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith($oldFull)) {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph with text '$oldFull'"
}

$pStart = $target.Range.Start

# ---- sanity check the exact original run text we are about to touch ----
$checkRange = $d.Range($pStart, $pStart + $oldFull.Length)
if ($checkRange.Text -ne $oldFull) {
    throw "Unexpected paragraph content: [$($checkRange.Text)]"
}

# ---- Step 1: rewrite run 1 ----
# "Product Owner<nbsp>: "  ->  "Produkt Owner<nbsp>: Herr"
$run1Range = $d.Range($pStart, $pStart + $oldRun1.Length)
$newRun1 = "Produkt Owner" + $nbsp + ": Herr"
$run1Range.Text = $newRun1

# ---- Step 2: rewrite run 2 (add the leading space before "Prof.") ----
$remStart = $pStart + $newRun1.Length
$remRange = $d.Range($remStart, $remStart + $oldRun2.Length)
if ($remRange.Text -ne $oldRun2) {
    throw "Unexpected remainder content: [$($remRange.Text)]"
}
$newRemainder = " " + $oldRun2
$remRange.Text = $newRemainder

# ---- Step 3: split the rewritten remainder into four runs ----
#   " Prof. Dr." | " " | "Christian " | "Kruse"
$splitOffsets = @(10, 11, 21, $newRemainder.Length)
$prevOffset = 0
foreach ($offset in $splitOffsets) {
    $a = $remStart + $prevOffset
    $b = $remStart + $offset
    $piece = $d.Range($a, $b)
    # Toggling a character-format property and restoring it forces Word
    # to keep this span as its own run without changing its appearance.
    $piece.Font.Bold = 1
    $piece.Font.Bold = 0
    $prevOffset = $offset
}

Write-Host "Result: [$($target.Range.Text)]"
